$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1686.9445
$ws.Range("J17").Value = 1824.5
$ws.Range("L17").Value = 5473.5
$ws.Range("N17").Value = -5809.5
$ws.Range("H86").Value = 1729.1538
$ws.Range("I86").Value = 1643.5454
$ws.Range("K86").Value = 1643.5454
$ws.Range("M86").Value = -520.5454
$ws.Range("H89").Value = 1729.1538
$ws.Range("I89").Value = 1643.5454
$ws.Range("K89").Value = 8217.726999999999
$ws.Range("M89").Value = -2601.726999999999
$ws.Range("H97").Value = 1995.1
$ws.Range("I97").Value = 989.5
$ws.Range("J97").Value = 2246.5
$ws.Range("K97").Value = 2968.5
$ws.Range("L97").Value = 6739.5
$ws.Range("M97").Value = -2472.5
$ws.Range("N97").Value = -7731.5
$ws.Range("H112").Value = 2616.6047
$ws.Range("J112").Value = 2359.6052
$ws.Range("L112").Value = 7078.8156
$ws.Range("N112").Value = -9294.8156
$ws.Range("H118").Value = 248
$ws.Range("I118").Value = 248
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 744
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = ""
$ws.Range("N118").Value = 913
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1685856.2
$ws.Range("I32").Value = 2390.7385
$ws.Range("J32").Value = 111111110
$ws.Range("K32").Value = 2390.7385
$ws.Range("L32").Value = 111111110
$ws.Range("M32").Value = -2103.7385
$ws.Range("N32").Value = -111111684
$ws.Range("H74").Value = 4518.8945
$ws.Range("I74").Value = 5011.357
$ws.Range("K74").Value = 5011.357
$ws.Range("M74").Value = -4137.357
$ws.Range("H77").Value = 4518.8945
$ws.Range("I77").Value = 5011.357
$ws.Range("K77").Value = 25056.785
$ws.Range("M77").Value = -20688.785
$ws.Range("H97").Value = 818.05884
$ws.Range("I97").Value = 327.13333
$ws.Range("J97").Value = 4500
$ws.Range("K97").Value = 327.13333
$ws.Range("L97").Value = 4500
$ws.Range("M97").Value = 168.86667
$ws.Range("N97").Value = -5492
$ws.Range("H102").Value = 1714.7778
$ws.Range("I102").Value = 1750.7693
$ws.Range("J102").Value = 779
$ws.Range("K102").Value = 1750.7693
$ws.Range("L102").Value = 779
$ws.Range("M102").Value = -128.7692999999999
$ws.Range("N102").Value = -4023
$ws.Range("H122").Value = 2647.3333
$ws.Range("I122").Value = 1881.85
$ws.Range("K122").Value = 5645.549999999999
$ws.Range("M122").Value = -3195.549999999999
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3367.5386
$ws.Range("I86").Value = 1465.0952
$ws.Range("K86").Value = 1465.0952
$ws.Range("M86").Value = -342.0952
$ws.Range("H89").Value = 3367.5386
$ws.Range("I89").Value = 1465.0952
$ws.Range("K89").Value = 7325.476
$ws.Range("M89").Value = -1709.476
$ws.Range("H94").Value = 3153.5833
$ws.Range("J94").Value = 7019.4614
$ws.Range("L94").Value = 7019.4614
$ws.Range("N94").Value = -7921.4614
$ws.Range("H105").Value = 1684.875
$ws.Range("I105").Value = 1671.174
$ws.Range("K105").Value = 1671.174
$ws.Range("M105").Value = 75.82600000000002
$ws.Range("H107").Value = 7148948.5
$ws.Range("I107").Value = 7698098.5
$ws.Range("K107").Value = 7698098.5
$ws.Range("M107").Value = -7696178.5
$ws.Range("H134").Value = 1357899.5
$ws.Range("I134").Value = 1857938.9
$ws.Range("K134").Value = 5573816.699999999
$ws.Range("M134").Value = -5571281.699999999
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1637
$ws.Range("I10").Value = 1332.6666
$ws.Range("J10").Value = 2550
$ws.Range("K10").Value = 1332.6666
$ws.Range("L10").Value = 2550
$ws.Range("M10").Value = -1193.6666
$ws.Range("N10").Value = -2828
$ws.Range("H62").Value = 11917.167
$ws.Range("J62").Value = 5499.5
$ws.Range("L62").Value = 5499.5
$ws.Range("N62").Value = -6747.5
$ws.Range("H65").Value = 11917.167
$ws.Range("J65").Value = 5499.5
$ws.Range("L65").Value = 27497.5
$ws.Range("N65").Value = -33737.5
$ws.Range("H105").Value = 37038944
$ws.Range("I105").Value = 47620924
$ws.Range("J105").Value = 2003.5
$ws.Range("K105").Value = 47620924
$ws.Range("L105").Value = 2003.5
$ws.Range("M105").Value = -47619177
$ws.Range("N105").Value = -5497.5
$ws.Range("H132").Value = 19165.416
$ws.Range("I132").Value = 13877.75
$ws.Range("J132").Value = 21809.25
$ws.Range("K132").Value = 41633.25
$ws.Range("L132").Value = 65427.75
$ws.Range("M132").Value = -39103.25
$ws.Range("N132").Value = -70487.75
$ws.Range("H134").Value = 71435910
$ws.Range("I134").Value = 90914370
$ws.Range("J134").Value = 14883.333
$ws.Range("K134").Value = 272743110
$ws.Range("L134").Value = 44649.999
$ws.Range("M134").Value = -272740575
$ws.Range("N134").Value = -49719.999
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 14362.615
$ws.Range("I64").Value = 9913
$ws.Range("K64").Value = 29739
$ws.Range("M64").Value = -29469
$ws.Range("H67").Value = 14362.615
$ws.Range("I67").Value = 9913
$ws.Range("K67").Value = 29739
$ws.Range("M67").Value = -28803
$ws.Range("H105").Value = 26833.166
$ws.Range("J105").Value = 26833.166
$ws.Range("L105").Value = 80499.49800000001
$ws.Range("N105").Value = -85741.49800000001
$ws.Range("H134").Value = 71449250
$ws.Range("I134").Value = 76942850
$ws.Range("K134").Value = 230828550
$ws.Range("M134").Value = -230823480
$ws.Range("H139").Value = 12502428
$ws.Range("J139").Value = 4299
$ws.Range("L139").Value = 12897
$ws.Range("N139").Value = -23177
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5022.727
$ws.Range("I80").Value = 3550
$ws.Range("J80").Value = 19750
$ws.Range("K80").Value = 3550
$ws.Range("L80").Value = 19750
$ws.Range("M80").Value = -2552
$ws.Range("N80").Value = -21746
$ws.Range("H83").Value = 5022.727
$ws.Range("I83").Value = 3550
$ws.Range("J83").Value = 19750
$ws.Range("K83").Value = 17750
$ws.Range("L83").Value = 98750
$ws.Range("M83").Value = -12758
$ws.Range("N83").Value = -108734
$ws.Range("H122").Value = 4571.2915
$ws.Range("I122").Value = 3262.1875
$ws.Range("K122").Value = 9786.5625
$ws.Range("M122").Value = -7336.5625
$ws.Range("H132").Value = 47623620
$ws.Range("I132").Value = 71432540
$ws.Range("K132").Value = 214297620
$ws.Range("M132").Value = -214295090
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 17110
$ws.Range("I2").Value = 17110
$ws.Range("K2").Value = 17110
$ws.Range("M2").Value = -16998
$ws.Range("H40").Value = 8124.909
$ws.Range("I40").Value = 8137.5
$ws.Range("J40").Value = 7999
$ws.Range("K40").Value = 8137.5
$ws.Range("L40").Value = 7999
$ws.Range("M40").Value = -8001.5
$ws.Range("N40").Value = -8271
$ws.Range("H82").Value = 3562.5
$ws.Range("I82").Value = 2090.9092
$ws.Range("K82").Value = 2090.9092
$ws.Range("M82").Value = -1729.9092
$ws.Range("H85").Value = 3562.5
$ws.Range("I85").Value = 2090.9092
$ws.Range("K85").Value = 2090.9092
$ws.Range("M85").Value = -842.9092000000001
$ws.Range("H93").Value = 1396
$ws.Range("I93").Value = 1304.0435
$ws.Range("J93").Value = 1924.75
$ws.Range("K93").Value = 1304.0435
$ws.Range("L93").Value = 1924.75
$ws.Range("M93").Value = -56.04349999999999
$ws.Range("N93").Value = -4420.75
$ws.Range("H100").Value = 3144.65
$ws.Range("I100").Value = 4862
$ws.Range("J100").Value = 1999.75
$ws.Range("K100").Value = 4862
$ws.Range("L100").Value = 1999.75
$ws.Range("M100").Value = -4321
$ws.Range("N100").Value = -3081.75
$ws.Range("H132").Value = 6132.5557
$ws.Range("I132").Value = 4448.75
$ws.Range("J132").Value = 7479.6
$ws.Range("K132").Value = 13346.25
$ws.Range("L132").Value = 22438.8
$ws.Range("M132").Value = -10816.25
$ws.Range("N132").Value = -27498.8
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -94060
$ws.Range("H136").Value = 58832056
$ws.Range("I136").Value = 142867260
$ws.Range("K136").Value = 428601780
$ws.Range("M136").Value = -428599230
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16135897
$ws.Range("I136").Value = 17858646
$ws.Range("K136").Value = 53575938
$ws.Range("M136").Value = -53573388
